$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Emanuele Miorandi"
$ws.Range("B10").Value = "Samuele Kettamier | SBARX"
$ws.Range("C10").Value = "Alessio Zandonai | SBARX"
$ws.Range("D10").Value = "Luca Frasca | Clitoriders"
$ws.Range("E10").Value = "Federico Andreis | iMontagna"
$ws.Range("F10").Value = "Danny Giordani | I Magnifici"
